$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "23×16="
$t.Cell(1, 2).Range.Text = "14×64="
$t.Cell(1, 3).Range.Text = "66×89="
$t.Cell(1, 4).Range.Text = "99×74="
$t.Cell(1, 5).Range.Text = "84×65="
$t.Cell(5, 1).Range.Text = "45×98="
$t.Cell(5, 2).Range.Text = "55×65="
$t.Cell(5, 3).Range.Text = "95×38="
$t.Cell(5, 4).Range.Text = "13×98="
$t.Cell(5, 5).Range.Text = "24×68="
$t.Cell(10, 1).Range.Text = "14×50="
$t.Cell(10, 2).Range.Text = "93×87="
$t.Cell(10, 3).Range.Text = "63×73="
$t.Cell(10, 4).Range.Text = "62×44="
$t.Cell(10, 5).Range.Text = "62×89="
$t.Cell(15, 1).Range.Text = "80×73="
$t.Cell(15, 2).Range.Text = "34×37="
$t.Cell(15, 3).Range.Text = "82×44="
$t.Cell(15, 4).Range.Text = "95×89="
$t.Cell(15, 5).Range.Text = "25×54="
$t.Cell(20, 1).Range.Text = "68×93="
$t.Cell(20, 2).Range.Text = "94×23="
$t.Cell(20, 3).Range.Text = "98×50="
$t.Cell(20, 4).Range.Text = "80×44="
$t.Cell(20, 5).Range.Text = "92×74="
